$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A and append the next day's data below it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($lastRow, 1).Value2 + 1
$ws.Cells.Item($newRow, 2).Value = 120
$ws.Cells.Item($newRow, 3).Value = 132
$ws.Cells.Item($newRow, 4).Value = 122

$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
